# Auto-generated Excel COM-interop script
# Applies the Gilgamesh_Profits.xlsx market-data refresh described in the commit diff.
# For each affected leve row, updates currentAveragePrice* / Leve* columns (H:N) to the
# latest scraped market values; some rows gain or lose a profit cell (M/N) depending on
# whether that computed value is present in the refreshed data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1599.25
$ws.Cells.Item(19, 9).Value = 1498
$ws.Cells.Item(19, 11).Value = 1498
$ws.Cells.Item(19, 13).Value = -1323
$ws.Cells.Item(21, 8).Value = 29700
$ws.Cells.Item(21, 10).Value = 29700
$ws.Cells.Item(21, 12).Value = 29700
$ws.Cells.Item(21, 14).Value = -30636
$ws.Cells.Item(23, 8).Value = 29700
$ws.Cells.Item(23, 10).Value = 29700
$ws.Cells.Item(23, 12).Value = 29700
$ws.Cells.Item(23, 14).Value = -30168
$ws.Cells.Item(51, 8).Value = 8944.75
$ws.Cells.Item(51, 10).Value = 8593
$ws.Cells.Item(51, 12).Value = 8593
$ws.Cells.Item(51, 14).Value = -9561
$ws.Cells.Item(62, 8).Value = 3899.5
$ws.Cells.Item(62, 9).Value = 3899.5
$ws.Cells.Item(62, 11).Value = 3899.5
$ws.Cells.Item(62, 13).Value = -3275.5
$ws.Cells.Item(65, 8).Value = 3899.5
$ws.Cells.Item(65, 9).Value = 3899.5
$ws.Cells.Item(65, 11).Value = 19497.5
$ws.Cells.Item(65, 13).Value = -16377.5
$ws.Cells.Item(70, 8).Value = 2432.8333
$ws.Cells.Item(70, 10).Value = 2000
$ws.Cells.Item(70, 12).Value = 6000
$ws.Cells.Item(70, 14).Value = -6540
$ws.Cells.Item(73, 8).Value = 2432.8333
$ws.Cells.Item(73, 10).Value = 2000
$ws.Cells.Item(73, 12).Value = 6000
$ws.Cells.Item(73, 14).Value = -7872
$ws.Cells.Item(86, 8).Value = 333333660
$ws.Cells.Item(86, 9).Value = 333333660
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 333333660
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = -333332537
$ws.Cells.Item(86, 14).ClearContents()
$ws.Cells.Item(87, 8).Value = 133333.33
$ws.Cells.Item(87, 10).Value = 133333.33
$ws.Cells.Item(87, 12).Value = 133333.33
$ws.Cells.Item(87, 14).Value = -135829.33
$ws.Cells.Item(89, 8).Value = 333333660
$ws.Cells.Item(89, 9).Value = 333333660
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 1666668300
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = -1666662684
$ws.Cells.Item(89, 14).ClearContents()
$ws.Cells.Item(90, 8).Value = 133333.33
$ws.Cells.Item(90, 10).Value = 133333.33
$ws.Cells.Item(90, 12).Value = 399999.99
$ws.Cells.Item(90, 14).Value = -412479.99
$ws.Cells.Item(98, 8).Value = 1886.1351
$ws.Cells.Item(98, 9).Value = 1888.5278
$ws.Cells.Item(98, 11).Value = 1888.5278
$ws.Cells.Item(98, 13).Value = -390.5278000000001
$ws.Cells.Item(106, 8).Value = 1888.4117
$ws.Cells.Item(106, 9).Value = 1726.8667
$ws.Cells.Item(106, 11).Value = 1726.8667
$ws.Cells.Item(106, 13).Value = -1095.8667
$ws.Cells.Item(113, 8).Value = 3000.9
$ws.Cells.Item(113, 9).Value = 2127.5
$ws.Cells.Item(113, 11).Value = 2127.5
$ws.Cells.Item(113, 13).Value = 1126.5
$ws.Cells.Item(116, 8).Value = 5150
$ws.Cells.Item(116, 9).Value = 2800
$ws.Cells.Item(116, 10).Value = 7500
$ws.Cells.Item(116, 11).Value = 2800
$ws.Cells.Item(116, 12).Value = 7500
$ws.Cells.Item(116, 13).Value = 642
$ws.Cells.Item(116, 14).Value = -14384
$ws.Cells.Item(122, 8).Value = 1886.1351
$ws.Cells.Item(122, 9).Value = 1888.5278
$ws.Cells.Item(122, 11).Value = 5665.5834
$ws.Cells.Item(122, 13).Value = -3215.5834
$ws.Cells.Item(137, 8).Value = 22786.34
$ws.Cells.Item(137, 9).Value = 26336.385
$ws.Cells.Item(137, 11).Value = 79009.155
$ws.Cells.Item(137, 13).Value = -76459.155
$ws.Cells.Item(141, 8).Value = 4987.75
$ws.Cells.Item(141, 9).Value = 3137.3044
$ws.Cells.Item(141, 10).Value = 13499.8
$ws.Cells.Item(141, 11).Value = 9411.913199999999
$ws.Cells.Item(141, 12).Value = 40499.39999999999
$ws.Cells.Item(141, 13).Value = -4231.913199999999
$ws.Cells.Item(141, 14).Value = -50859.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1422.8572
$ws.Cells.Item(2, 9).Value = 1453.6666
$ws.Cells.Item(2, 10).Value = 1399.75
$ws.Cells.Item(2, 11).Value = 1453.6666
$ws.Cells.Item(2, 12).Value = 1399.75
$ws.Cells.Item(2, 13).Value = -1340.6666
$ws.Cells.Item(2, 14).Value = -1625.75
$ws.Cells.Item(32, 8).Value = 2940.0576
$ws.Cells.Item(32, 9).Value = 2940.0576
$ws.Cells.Item(32, 11).Value = 2940.0576
$ws.Cells.Item(32, 13).Value = -2653.0576
$ws.Cells.Item(74, 8).Value = 244085.66
$ws.Cells.Item(74, 9).Value = 464232.34
$ws.Cells.Item(74, 10).Value = 3925.6365
$ws.Cells.Item(74, 11).Value = 464232.34
$ws.Cells.Item(74, 12).Value = 3925.6365
$ws.Cells.Item(74, 13).Value = -463358.34
$ws.Cells.Item(74, 14).Value = -5673.636500000001
$ws.Cells.Item(77, 8).Value = 244085.66
$ws.Cells.Item(77, 9).Value = 464232.34
$ws.Cells.Item(77, 10).Value = 3925.6365
$ws.Cells.Item(77, 11).Value = 2321161.7
$ws.Cells.Item(77, 12).Value = 19628.1825
$ws.Cells.Item(77, 13).Value = -2316793.7
$ws.Cells.Item(77, 14).Value = -28364.1825
$ws.Cells.Item(110, 8).Value = 2650.6296
$ws.Cells.Item(110, 9).Value = 1564.8
$ws.Cells.Item(110, 11).Value = 1564.8
$ws.Cells.Item(110, 13).Value = 480.2
$ws.Cells.Item(116, 8).Value = 1422.8572
$ws.Cells.Item(116, 9).Value = 1453.6666
$ws.Cells.Item(116, 10).Value = 1399.75
$ws.Cells.Item(116, 11).Value = 1453.6666
$ws.Cells.Item(116, 12).Value = 1399.75
$ws.Cells.Item(116, 13).Value = 840.3334
$ws.Cells.Item(116, 14).Value = -5987.75
$ws.Cells.Item(122, 8).Value = 4509.1875
$ws.Cells.Item(122, 9).Value = 4022.64
$ws.Cells.Item(122, 10).Value = 6246.857
$ws.Cells.Item(122, 11).Value = 12067.92
$ws.Cells.Item(122, 12).Value = 18740.571
$ws.Cells.Item(122, 13).Value = -9617.92
$ws.Cells.Item(122, 14).Value = -23640.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1422.8572
$ws.Cells.Item(3, 9).Value = 1453.6666
$ws.Cells.Item(3, 10).Value = 1399.75
$ws.Cells.Item(3, 11).Value = 1453.6666
$ws.Cells.Item(3, 12).Value = 1399.75
$ws.Cells.Item(3, 13).Value = -1339.6666
$ws.Cells.Item(3, 14).Value = -1627.75
$ws.Cells.Item(22, 8).Value = 1098.8572
$ws.Cells.Item(22, 9).Value = 1032
$ws.Cells.Item(22, 10).Value = 1500
$ws.Cells.Item(22, 11).Value = 1032
$ws.Cells.Item(22, 12).Value = 1500
$ws.Cells.Item(22, 13).Value = -859
$ws.Cells.Item(22, 14).Value = -1846
$ws.Cells.Item(86, 8).Value = 3574.5
$ws.Cells.Item(86, 9).Value = 3449.3333
$ws.Cells.Item(86, 10).Value = 3950
$ws.Cells.Item(86, 11).Value = 3449.3333
$ws.Cells.Item(86, 12).Value = 3950
$ws.Cells.Item(86, 13).Value = -2326.3333
$ws.Cells.Item(86, 14).Value = -6196
$ws.Cells.Item(89, 8).Value = 3574.5
$ws.Cells.Item(89, 9).Value = 3449.3333
$ws.Cells.Item(89, 10).Value = 3950
$ws.Cells.Item(89, 11).Value = 17246.6665
$ws.Cells.Item(89, 12).Value = 19750
$ws.Cells.Item(89, 13).Value = -11630.6665
$ws.Cells.Item(89, 14).Value = -30982
$ws.Cells.Item(105, 8).Value = 14447276
$ws.Cells.Item(105, 9).Value = 835626.75
$ws.Cells.Item(105, 10).Value = 41670576
$ws.Cells.Item(105, 11).Value = 835626.75
$ws.Cells.Item(105, 12).Value = 41670576
$ws.Cells.Item(105, 13).Value = -833879.75
$ws.Cells.Item(105, 14).Value = -41674070
$ws.Cells.Item(134, 8).Value = 6539.8374
$ws.Cells.Item(134, 9).Value = 6263.8
$ws.Cells.Item(134, 11).Value = 18791.4
$ws.Cells.Item(134, 13).Value = -16256.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5531.657
$ws.Cells.Item(31, 9).Value = 4275.4375
$ws.Cells.Item(31, 11).Value = 4275.4375
$ws.Cells.Item(31, 13).Value = -3980.4375
$ws.Cells.Item(34, 8).Value = 5531.657
$ws.Cells.Item(34, 9).Value = 4275.4375
$ws.Cells.Item(34, 11).Value = 4275.4375
$ws.Cells.Item(34, 13).Value = -4073.4375
$ws.Cells.Item(100, 8).Value = 117890
$ws.Cells.Item(100, 10).Value = 117890
$ws.Cells.Item(100, 12).Value = 117890
$ws.Cells.Item(100, 14).Value = -120054
$ws.Cells.Item(105, 8).Value = 1705.8462
$ws.Cells.Item(105, 10).Value = 2675
$ws.Cells.Item(105, 12).Value = 2675
$ws.Cells.Item(105, 14).Value = -6169
$ws.Cells.Item(107, 8).Value = 656.8333
$ws.Cells.Item(107, 9).Value = 656.8333
$ws.Cells.Item(107, 11).Value = 656.8333
$ws.Cells.Item(107, 13).Value = 1263.1667
$ws.Cells.Item(122, 8).Value = 3280.5386
$ws.Cells.Item(122, 9).Value = 1801.4117
$ws.Cells.Item(122, 10).Value = 6074.4443
$ws.Cells.Item(122, 11).Value = 5404.2351
$ws.Cells.Item(122, 12).Value = 18223.3329
$ws.Cells.Item(122, 13).Value = -2954.2351
$ws.Cells.Item(122, 14).Value = -23123.3329
$ws.Cells.Item(132, 8).Value = 2351.825
$ws.Cells.Item(132, 9).Value = 2127.9062
$ws.Cells.Item(132, 10).Value = 3247.5
$ws.Cells.Item(132, 11).Value = 6383.7186
$ws.Cells.Item(132, 12).Value = 9742.5
$ws.Cells.Item(132, 13).Value = -3853.7186
$ws.Cells.Item(132, 14).Value = -14802.5
$ws.Cells.Item(134, 8).Value = 3018.4443
$ws.Cells.Item(134, 10).Value = 5249
$ws.Cells.Item(134, 12).Value = 15747
$ws.Cells.Item(134, 14).Value = -20817

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 3486255.2
$ws.Cells.Item(4, 9).Value = 3103785.8
$ws.Cells.Item(4, 10).Value = 6240035.5
$ws.Cells.Item(4, 11).Value = 9311357.399999999
$ws.Cells.Item(4, 12).Value = 18720106.5
$ws.Cells.Item(4, 13).Value = -9311245.399999999
$ws.Cells.Item(4, 14).Value = -18720330.5
$ws.Cells.Item(5, 8).Value = 864.5
$ws.Cells.Item(5, 9).Value = 738.44446
$ws.Cells.Item(5, 10).Value = 1999
$ws.Cells.Item(5, 11).Value = 2215.33338
$ws.Cells.Item(5, 12).Value = 5997
$ws.Cells.Item(5, 13).Value = -2103.33338
$ws.Cells.Item(5, 14).Value = -6221
$ws.Cells.Item(95, 8).Value = 8622.5
$ws.Cells.Item(95, 10).Value = 8622.5
$ws.Cells.Item(95, 12).Value = 25867.5
$ws.Cells.Item(95, 14).Value = -29985.5
$ws.Cells.Item(114, 8).Value = 947.4
$ws.Cells.Item(114, 9).Value = 187.5
$ws.Cells.Item(114, 10).Value = 2087.25
$ws.Cells.Item(114, 11).Value = 562.5
$ws.Cells.Item(114, 12).Value = 6261.75
$ws.Cells.Item(114, 13).Value = 2691.5
$ws.Cells.Item(114, 14).Value = -12769.75
$ws.Cells.Item(125, 8).Value = 1999.5
$ws.Cells.Item(125, 9).Value = 1999.5
$ws.Cells.Item(125, 11).Value = 5998.5
$ws.Cells.Item(125, 13).Value = -1078.5
$ws.Cells.Item(131, 8).Value = 22312.428
$ws.Cells.Item(131, 9).Value = 37225
$ws.Cells.Item(131, 11).Value = 111675
$ws.Cells.Item(131, 13).Value = -106635
$ws.Cells.Item(135, 8).Value = 864.5
$ws.Cells.Item(135, 9).Value = 738.44446
$ws.Cells.Item(135, 10).Value = 1999
$ws.Cells.Item(135, 11).Value = 6646.00014
$ws.Cells.Item(135, 12).Value = 17991
$ws.Cells.Item(135, 13).Value = -4111.00014
$ws.Cells.Item(135, 14).Value = -23061
$ws.Cells.Item(138, 8).Value = 3291.7
$ws.Cells.Item(138, 10).Value = 7998.5
$ws.Cells.Item(138, 12).Value = 23995.5
$ws.Cells.Item(138, 14).Value = -34275.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 111114056
$ws.Cells.Item(80, 9).Value = 125002500
$ws.Cells.Item(80, 11).Value = 125002500
$ws.Cells.Item(80, 13).Value = -125001502
$ws.Cells.Item(83, 8).Value = 111114056
$ws.Cells.Item(83, 9).Value = 125002500
$ws.Cells.Item(83, 11).Value = 625012500
$ws.Cells.Item(83, 13).Value = -625007508
$ws.Cells.Item(107, 8).Value = 1002
$ws.Cells.Item(107, 9).Value = 1002
$ws.Cells.Item(107, 11).Value = 1002
$ws.Cells.Item(107, 13).Value = 918
$ws.Cells.Item(122, 8).Value = 6159.0586
$ws.Cells.Item(122, 9).Value = 4882.4546
$ws.Cells.Item(122, 11).Value = 14647.3638
$ws.Cells.Item(122, 13).Value = -12197.3638
$ws.Cells.Item(126, 8).Value = 7689.1113
$ws.Cells.Item(126, 9).Value = 2535.5
$ws.Cells.Item(126, 11).Value = 7606.5
$ws.Cells.Item(126, 13).Value = -5136.5
$ws.Cells.Item(136, 8).Value = 9931.786
$ws.Cells.Item(136, 10).Value = 9931.786
$ws.Cells.Item(136, 12).Value = 29795.358
$ws.Cells.Item(136, 14).Value = -34895.358

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(56, 8).Value = 29435.715
$ws.Cells.Item(56, 9).Value = 22410
$ws.Cells.Item(56, 11).Value = 22410
$ws.Cells.Item(56, 13).Value = -21719
$ws.Cells.Item(123, 8).Value = 20390
$ws.Cells.Item(123, 9).Value = 20390
$ws.Cells.Item(123, 11).Value = 20390
$ws.Cells.Item(123, 13).Value = -15490
$ws.Cells.Item(136, 8).Value = 5243.4287
$ws.Cells.Item(136, 10).Value = 4282.8335
$ws.Cells.Item(136, 12).Value = 12848.5005
$ws.Cells.Item(136, 14).Value = -17948.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(81, 8).Value = 3427.2856
$ws.Cells.Item(81, 9).Value = 1177.6
$ws.Cells.Item(81, 10).Value = 6023.077
$ws.Cells.Item(81, 11).Value = 2355.2
$ws.Cells.Item(81, 12).Value = 12046.154
$ws.Cells.Item(81, 13).Value = -1294.2
$ws.Cells.Item(81, 14).Value = -14168.154
$ws.Cells.Item(84, 8).Value = 3427.2856
$ws.Cells.Item(84, 9).Value = 1177.6
$ws.Cells.Item(84, 10).Value = 6023.077
$ws.Cells.Item(84, 11).Value = 11776
$ws.Cells.Item(84, 12).Value = 60230.77
$ws.Cells.Item(84, 13).Value = -6472
$ws.Cells.Item(84, 14).Value = -70838.77
$ws.Cells.Item(122, 8).Value = 8930643
$ws.Cells.Item(122, 9).Value = 2075.6316
$ws.Cells.Item(122, 11).Value = 6226.8948
$ws.Cells.Item(122, 13).Value = -3776.8948

